$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column values that look numeric stay as plain text (matching source data)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "98.936.92"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "3.288.96"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "254.37"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "625.74"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "1.45"
$ws.Range("E7").Value = "  +22.63%  "
$ws.Range("E8").Value = "  +4.26%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "0.977"
$ws.Range("E10").Value = "  +22.72%  "
$ws.Range("D11").Value = "3.288.88"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").Value = "0.204"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "40.60"
$ws.Range("E13").Value = "  +13.42%  "
$ws.Range("D14").Value = "98.717.61"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "0.0000249"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "3.906.34"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "3.288.49"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").Value = "15.56"
$ws.Range("E20").Value = "  +5.28%  "
$ws.Range("D21").Value = "6.40"
$ws.Range("E21").Value = "  +9.10%  "
$ws.Range("D22").Value = "488.77"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").Value = "9.37"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "0.345"
$ws.Range("E25").Value = "  +39.02%  "
$ws.Range("D26").Value = "5.69"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "12.13"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Value = "3.459.69"
$ws.Range("E29").Value = "  -3.12%  "
$ws.Range("D30").Value = "0.146"
$ws.Range("E30").Value = "  +19.03%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "0.190"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").Value = "10.63"
$ws.Range("E33").Value = "  +15.45%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "27.93"
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("D36").Value = "0.480"
$ws.Range("E36").Value = "  +7.17%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "7.29"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "493.59"
$ws.Range("E40").Value = "  -5.72%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "3.90"
$ws.Range("E42").Value = "  +8.06%  "
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D45").Value = "0.779"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "3.13"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("D47").Value = "159.06"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").Value = "4.82"
$ws.Range("E49").Value = "  +7.17%  "
$ws.Range("E50").Value = "  +15.48%  "
$ws.Range("D51").Value = "0.849"
$ws.Range("E51").Value = "  +7.40%  "
